$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsAsset    = $wb.Worksheets.Item("Asset_Cal_Info")

# ---------------------------------------------------------------------------
# Asset_Cal_Info sheet: instrument reference designator id correction
# (CP05MOAS-GL002 -> CP05MOAS-GL340) + deployment number correction,
# across every populated data row (Ref Des column A, Deployment Number column C).
# ---------------------------------------------------------------------------
$wsAsset.Range("A2").Value  = "CP05MOAS-GL340-01-ADCPAM000"
$wsAsset.Range("A3").Value  = "CP05MOAS-GL340-01-ADCPAM000"
$wsAsset.Range("A4").Value  = "CP05MOAS-GL340-01-ADCPAM000"
$wsAsset.Range("A5").Value  = "CP05MOAS-GL340-01-ADCPAM000"

$wsAsset.Range("A7").Value  = "CP05MOAS-GL340-02-FLORTM000"
$wsAsset.Range("A8").Value  = "CP05MOAS-GL340-02-FLORTM000"
$wsAsset.Range("A9").Value  = "CP05MOAS-GL340-02-FLORTM000"
$wsAsset.Range("A10").Value = "CP05MOAS-GL340-02-FLORTM000"

$wsAsset.Range("A12").Value = "CP05MOAS-GL340-03-CTDGVM000"
$wsAsset.Range("A14").Value = "CP05MOAS-GL340-04-DOSTAM000"
$wsAsset.Range("A16").Value = "CP05MOAS-GL340-05-PARADM000"
$wsAsset.Range("A18").Value = "CP05MOAS-GL340-00-ENG000000"

$wsAsset.Range("C2").Value  = 1
$wsAsset.Range("C3").Value  = 1
$wsAsset.Range("C4").Value  = 1
$wsAsset.Range("C5").Value  = 1
$wsAsset.Range("C7").Value  = 1
$wsAsset.Range("C8").Value  = 1
$wsAsset.Range("C9").Value  = 1
$wsAsset.Range("C10").Value = 1
$wsAsset.Range("C12").Value = 1
$wsAsset.Range("C14").Value = 1
$wsAsset.Range("C16").Value = 1
$wsAsset.Range("C18").Value = 1

# ---------------------------------------------------------------------------
# Moorings sheet: the mooring was re-deployed under a new glider id
# (CP05MOAS-GL002 -> CP05MOAS-GL340) on deployment 1 (was 2).
# ---------------------------------------------------------------------------
$wsMoorings.Range("A2").Value = "CP05MOAS-GL340"
$wsMoorings.Range("C2").Value = 1

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: Moorings becomes the active tab,
# with a fresh selection on each sheet.
# ---------------------------------------------------------------------------
$wsAsset.Range("A40").Select()
$wsMoorings.Activate()
$wsMoorings.Range("E12").Select()
